$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.7287194209349384
$ws.Range("C2").Value = 10990084.13351303
$ws.Range("D2").Value = 16.98373111632243
$ws.Range("E2").Value = 71517.89157740913
$ws.Range("G2").Value = 11061619.73754098
